# Generate Report for Archive
#
# Localization status moved from "Ready for handoff" to "In Translation"
# for the zh-cn / de-de targets. The same shared string is used by the
# "Overview" sheet's per-language status columns (E = zh-cn, F = de-de)
# and by the "Status" column (C) on each language detail sheet, so all
# four cells are updated together.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# The shorter status text lets Excel shrink the (previously autosized)
# status columns back down.
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
